$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.530.35"
$ws.Range("E2").Value = "  -2.85%  "

$ws.Range("D3").Value = "1.670.34"
$ws.Range("E3").Value = "  -2.23%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.42%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.96%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5150"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.96%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.007"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.37%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.06458"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2569"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.25%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.99"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.40%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07663"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.37%  "

$ws.Range("D12").Value = "1.691.36"
$ws.Range("E12").Value = "  -1.25%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.343"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.20%  "

$ws.Range("D14").Value = "1.899.40"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5564"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.05%  "

$ws.Range("D16").Value = "0.0₅8020"
$ws.Range("E16").Value = "  -2.12%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.76"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.19%  "

$ws.Range("D18").Value = "26.583.67"
$ws.Range("E18").Value = "  -2.64%  "

$ws.Range("E19").Value = "  +0.35%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "210.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.444"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.98%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.888"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.40%  "

$ws.Range("E24").Value = "  +0.29%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.47%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.714"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.90%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1167"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.984"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.83%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.72"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.75%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05205"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.29%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.264"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.353"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.45%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.203"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.47%  "

$ws.Range("E34").Value = "  -4.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.760"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.13%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9256"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.34%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5717"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.33%  "

$ws.Range("D39").Value = "1.151.10"
$ws.Range("E39").Value = "  +10.44%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01590"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.63%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.006"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.31%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8356"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.78%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.645"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.75%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.94"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.17%  "

$ws.Range("D45").Value = "1.809.67"
$ws.Range("E45").Value = "  -2.22%  "

$ws.Range("D46").Value = "0.0₈113"
$ws.Range("E46").Value = "  -2.64%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4491"
$ws.Range("D47").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.006"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.23%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.918"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.47%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05136"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.02%  "

Write-Output "Applied updates to cryptos sheet"
